$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2 through 18
# from serial date 45224 (2023-10-25) to 45233 (2023-11-03)
for ($row = 2; $row -le 18; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45224) {
        $cell.Value = 45233
    }
}
